# Insert a new weekly price record before the current row 185 ("Fruta /
# hortaliza, semanal" update) for Agrícola del Norte S.A. de Arica -
# Mandarina. Inserting the row pushes the existing rows 185-201 down to
# 186-202, matching the new dimension A1:T202.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = 1
$ws.Range("B185").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C185").Value = "Arica y Parinacota"
$ws.Range("D185").Value = 45194
$ws.Range("E185").Value = 15
$ws.Range("F185").Value = "Fruta"
$ws.Range("G185").Value = 100102
$ws.Range("H185").Value = "Cítricos"
$ws.Range("I185").Value = 100102004
$ws.Range("J185").Value = "Mandarina"
$ws.Range("K185").Value = "Murcott"
$ws.Range("L185").Value = "Segunda"
$ws.Range("M185").Value = 300
$ws.Range("N185").Value = 15000
$ws.Range("O185").Value = 16000
$ws.Range("P185").Value = 15500
$ws.Range("Q185").Value = "$/caja 20 kilos"
$ws.Range("R185").Value = "Región de Coquimbo"
$ws.Range("S185").Value = 775
$ws.Range("T185").Value = 20
